# general_stats.xlsx - update backtest stats rows (4,5,6) with new values
# coming from the unemployment-rate strategy, and clear the now-removed
# "vendido" (short) row (row 7) data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("periodo entero") updated values
$ws.Range("B4").Value = 0.05
$ws.Range("C4").Value = 0.03
$ws.Range("D4").Value = 0.04
$ws.Range("E4").Value = 0.99
$ws.Range("F4").Value = 11.58
$ws.Range("G4").Value = -20.47
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.03
$ws.Range("J4").Value = 0.04
$ws.Range("K4").Value = 0.8
$ws.Range("L4").Value = 9.1
$ws.Range("M4").Value = -20.47
$ws.Range("N4").Value = 0.05
$ws.Range("O4").Value = 0.03
$ws.Range("P4").Value = 0.04
$ws.Range("Q4").Value = 0.99
$ws.Range("R4").Value = 11.58
$ws.Range("S4").Value = -20.47

# Row 5 ("comprado") updated values
$ws.Range("B5").Value = 0.06
$ws.Range("C5").Value = 0.04
$ws.Range("D5").Value = 0.04
$ws.Range("E5").Value = 0.87
$ws.Range("F5").Value = 9.1
$ws.Range("G5").Value = -20.47
$ws.Range("H5").Value = 0.06
$ws.Range("I5").Value = 0.04
$ws.Range("J5").Value = 0.04
$ws.Range("K5").Value = 0.87
$ws.Range("L5").Value = 9.1
$ws.Range("M5").Value = -20.47
$ws.Range("N5").Value = 0.06
$ws.Range("O5").Value = 0.04
$ws.Range("P5").Value = 0.04
$ws.Range("Q5").Value = 0.87
$ws.Range("R5").Value = 9.1
$ws.Range("S5").Value = -20.47

# Row 6 ("en efectivo") updated values (H6:M6 stay 0, untouched)
$ws.Range("B6").Value = -0.01
$ws.Range("C6").Value = -0.02
$ws.Range("D6").Value = -0.01
$ws.Range("E6").Value = 1.49
$ws.Range("F6").Value = 11.58
$ws.Range("G6").Value = -11.98
$ws.Range("N6").Value = -0.01
$ws.Range("O6").Value = -0.02
$ws.Range("P6").Value = -0.01
$ws.Range("Q6").Value = 1.49
$ws.Range("R6").Value = 11.58
$ws.Range("S6").Value = -11.98

# Row 7 ("vendido") no longer has data - blank out B7:S7 while keeping the
# cells present as empty text cells (not entirely removed).
$ws.Range("B7:S7").Value = "'"
$ws.Range("B7:S7").Style = "Normal"
